$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "C2" = 0.012308799
    "D2" = 0.0041131
    "G2" = 0.0003083
    "H2" = 0.0003954

    "A3" = 0.0263105
    "B3" = 0.0207453
    "C3" = 0.0051289
    "D3" = 0.0060277
    "E3" = 0.131474
    "F3" = 0.0994753
    "G3" = 0.0004857
    "H3" = 0.0006566
    "I3" = 0.0001846
    "J3" = 0.0004467

    "A4" = 0.0435865
    "B4" = 0.0206726
    "C4" = 0.0073527
    "D4" = 0.0019882
    "E4" = 0.0953628
    "F4" = 0.097741501
    "I4" = 0.0005318
    "J4" = 0.0001778

    "A5" = 0.0260132
    "B5" = 0.0249152
    "C5" = 0.0045501
    "D5" = 0.0019824
    "E5" = 0.1434911
    "F5" = 0.0970755
    "I5" = 0.0003046
    "J5" = 0.0003493
    "K5" = 0.0002488
    "L5" = 0.0001388

    "A6" = 0.0223178
    "B6" = 0.0432158
    "C6" = 0.0066547
    "D6" = 0.001899
    "E6" = 0.1159043
    "F6" = 0.1134566

    "C7" = 0.0060368
    "D7" = 0.0057689
    "E7" = 0.1125136
    "F7" = 0.1190275

    "C8" = 0.0074249
    "D8" = 0.0022617

    "C9" = 0.0046888

    "C10" = 0.0129891
    "D10" = 0.0029162

    "C11" = 0.0051052
    "D11" = 0.0017404
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    $cell.ClearFormats()
    $cell.Value = $changes[$addr]
}
